$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("C2").Value = 1.173973083496094
$ws.Range("E2").Value = 1563.547187408227
$ws.Range("F2").Value = 0.07413226043880848
$ws.Range("G2").Value = 0.05665253491457791
$ws.Range("H2").Value = 0.04835475435085952
$ws.Range("I2").Value = 0.04356527353342693
$ws.Range("J2").Value = 0.03931883653220122
$ws.Range("K2").Value = 0.03721419189091164
$ws.Range("L2").Value = 0.03515448336341394
$ws.Range("M2").Value = 0.03436593252461427
$ws.Range("N2").Value = 0.03373797580173417
$ws.Range("O2").Value = 0.03280077365886915
$ws.Range("P2").Value = 0.03222708804034035
$ws.Range("Q2").Value = 0.03187862165922341
$ws.Range("R2").Value = 0.03164011354019894
$ws.Range("S2").Value = 0.03117954519611916
$ws.Range("T2").Value = 0.03093850863814578
$ws.Range("U2").Value = 0.03081998543866039
$ws.Range("V2").Value = 0.03068385546372995
$ws.Range("W2").Value = 0.03053249438010132
$ws.Range("X2").Value = 0.03050532668905769
$ws.Range("Y2").Value = 0.03047850267852293

# Row 3
$ws.Range("C3").Value = 1.111999034881592
$ws.Range("E3").Value = 1532.955287202603
$ws.Range("F3").Value = 0.07947233859252406
$ws.Range("G3").Value = 0.06420178573102718
$ws.Range("H3").Value = 0.05078747664476856
$ws.Range("I3").Value = 0.04593257105825917
$ws.Range("J3").Value = 0.04275806102706295
$ws.Range("K3").Value = 0.03948619307127019
$ws.Range("L3").Value = 0.03741482683367984
$ws.Range("M3").Value = 0.03640752496621434
$ws.Range("N3").Value = 0.03483950736982036
$ws.Range("O3").Value = 0.03334303958550467
$ws.Range("P3").Value = 0.03227850033285153
$ws.Range("Q3").Value = 0.03171939123935075
$ws.Range("R3").Value = 0.03126537217525002
$ws.Range("S3").Value = 0.03086341206727682
$ws.Range("T3").Value = 0.03052695433013307
$ws.Range("U3").Value = 0.03041065458517972
$ws.Range("V3").Value = 0.03013034348423771
$ws.Range("W3").Value = 0.03005949569128873
$ws.Range("X3").Value = 0.02999706921209486
$ws.Range("Y3").Value = 0.02988216934118135

# Row 4
$ws.Range("C4").Value = 1.119001865386963
$ws.Range("E4").Value = 1534.82336978748
$ws.Range("F4").Value = 0.07766070316551799
$ws.Range("G4").Value = 0.05665119251833779
$ws.Range("H4").Value = 0.04912561126866217
$ws.Range("I4").Value = 0.04374698841124713
$ws.Range("J4").Value = 0.04006241427669662
$ws.Range("K4").Value = 0.03725453695531966
$ws.Range("L4").Value = 0.03582198000058514
$ws.Range("M4").Value = 0.03402770501765689
$ws.Range("N4").Value = 0.0334607379776943
$ws.Range("O4").Value = 0.03244672192785427
$ws.Range("P4").Value = 0.03193517214605673
$ws.Range("Q4").Value = 0.03142268560304064
$ws.Range("R4").Value = 0.03116759852780791
$ws.Range("S4").Value = 0.03083764816724955
$ws.Range("T4").Value = 0.03044638981175303
$ws.Range("U4").Value = 0.03033215963462093
$ws.Range("V4").Value = 0.03019518707484485
$ws.Range("W4").Value = 0.03004332584504309
$ws.Range("X4").Value = 0.02995321341626812
$ws.Range("Y4").Value = 0.02991858420638362

# Row 5
$ws.Range("C5").Value = 1.286028623580933
$ws.Range("E5").Value = 1559.876389005842
$ws.Range("F5").Value = 0.08062226718419341
$ws.Range("G5").Value = 0.0633537106116964
$ws.Range("H5").Value = 0.04885099048309342
$ws.Range("I5").Value = 0.0461462751543393
$ws.Range("J5").Value = 0.04155734818721813
$ws.Range("K5").Value = 0.03892009594781526
$ws.Range("L5").Value = 0.03737166290408434
$ws.Range("M5").Value = 0.035770579176684
$ws.Range("N5").Value = 0.03461791909707495
$ws.Range("O5").Value = 0.03340974734907779
$ws.Range("P5").Value = 0.03263078761518345
$ws.Range("Q5").Value = 0.03227940953829144
$ws.Range("R5").Value = 0.03191476994046372
$ws.Range("S5").Value = 0.03148428337404723
$ws.Range("T5").Value = 0.03125431785849882
$ws.Range("U5").Value = 0.03096990897421737
$ws.Range("V5").Value = 0.03067989980174966
$ws.Range("W5").Value = 0.0306025852791727
$ws.Range("X5").Value = 0.03048011964441564
$ws.Range("Y5").Value = 0.03040694715410997

# Row 6
$ws.Range("C6").Value = 1.205009698867798
$ws.Range("E6").Value = 1535.809331875953
$ws.Range("F6").Value = 0.07593223771527903
$ws.Range("G6").Value = 0.06230481861077591
$ws.Range("H6").Value = 0.04931639185210891
$ws.Range("I6").Value = 0.04574745426883278
$ws.Range("J6").Value = 0.04206125111771723
$ws.Range("K6").Value = 0.0406606025460587
$ws.Range("L6").Value = 0.03803777412273155
$ws.Range("M6").Value = 0.0358273046763253
$ws.Range("N6").Value = 0.03465561115603576
$ws.Range("O6").Value = 0.03312582697502559
$ws.Range("P6").Value = 0.03210338089486468
$ws.Range("Q6").Value = 0.03163888400740018
$ws.Range("R6").Value = 0.03109460867688568
$ws.Range("S6").Value = 0.03066645506783555
$ws.Range("T6").Value = 0.03033736260899277
$ws.Range("U6").Value = 0.03029357948821901
$ws.Range("V6").Value = 0.03015799508820699
$ws.Range("W6").Value = 0.03004708399141228
$ws.Range("X6").Value = 0.03000927871637843
$ws.Range("Y6").Value = 0.02993780374027198

# Row 7
$ws.Range("C7").Value = 1.264973640441895
$ws.Range("E7").Value = 1510.84219040336
$ws.Range("F7").Value = 0.07867579259095119
$ws.Range("G7").Value = 0.06272024326533665
$ws.Range("H7").Value = 0.05084214781277165
$ws.Range("I7").Value = 0.04677467330991568
$ws.Range("J7").Value = 0.04259023124793607
$ws.Range("K7").Value = 0.03849377373044899
$ws.Range("L7").Value = 0.03558293430125814
$ws.Range("M7").Value = 0.03387870098712459
$ws.Range("N7").Value = 0.03279165570677926
$ws.Range("O7").Value = 0.03181192382540502
$ws.Range("P7").Value = 0.03141428826649796
$ws.Range("Q7").Value = 0.0309986122625769
$ws.Range("R7").Value = 0.03057176781126823
$ws.Range("S7").Value = 0.0300382812994116
$ws.Range("T7").Value = 0.03001801131379205
$ws.Range("U7").Value = 0.02975853950624173
$ws.Range("V7").Value = 0.02968652714912619
$ws.Range("W7").Value = 0.02959877741035058
$ws.Range("X7").Value = 0.02947716910445757
$ws.Range("Y7").Value = 0.02945111482267758

# Row 8
$ws.Range("C8").Value = 1.117023229598999
$ws.Range("E8").Value = 1538.637242243713
$ws.Range("F8").Value = 0.07328839078205636
$ws.Range("G8").Value = 0.06149848644769747
$ws.Range("H8").Value = 0.05093753687523453
$ws.Range("I8").Value = 0.04466662858333086
$ws.Range("J8").Value = 0.04176336153414211
$ws.Range("K8").Value = 0.037626230230417
$ws.Range("L8").Value = 0.03722889232093515
$ws.Range("M8").Value = 0.03582012304593573
$ws.Range("N8").Value = 0.03422849625528666
$ws.Range("O8").Value = 0.03315109969708527
$ws.Range("P8").Value = 0.03215288287683542
$ws.Range("Q8").Value = 0.03178163848094227
$ws.Range("R8").Value = 0.03145420902721943
$ws.Range("S8").Value = 0.03103973969889052
$ws.Range("T8").Value = 0.03075157449647825
$ws.Range("U8").Value = 0.03055100333669223
$ws.Range("V8").Value = 0.03030857548635659
$ws.Range("W8").Value = 0.03015940785870153
$ws.Range("X8").Value = 0.03002503019763114
$ws.Range("Y8").Value = 0.02999292869870784

# Row 9
$ws.Range("C9").Value = 1.10099458694458
$ws.Range("E9").Value = 1567.487900034292
$ws.Range("F9").Value = 0.07194588475046097
$ws.Range("G9").Value = 0.06072631922046724
$ws.Range("H9").Value = 0.05175798554688724
$ws.Range("I9").Value = 0.04543041048694335
$ws.Range("J9").Value = 0.04176378041218366
$ws.Range("K9").Value = 0.03890348223401793
$ws.Range("L9").Value = 0.03655728964999343
$ws.Range("M9").Value = 0.03392571392047779
$ws.Range("N9").Value = 0.03358903549570231
$ws.Range("O9").Value = 0.03287601650786492
$ws.Range("P9").Value = 0.03230871120813613
$ws.Range("Q9").Value = 0.03184564758787532
$ws.Range("R9").Value = 0.03157907498213484
$ws.Range("S9").Value = 0.0313371440312802
$ws.Range("T9").Value = 0.03107813220950995
$ws.Range("U9").Value = 0.03094341213213484
$ws.Range("V9").Value = 0.03086118458976378
$ws.Range("W9").Value = 0.03068746042005874
$ws.Range("X9").Value = 0.03061174353015877
$ws.Range("Y9").Value = 0.03055531968877761

# Row 10
$ws.Range("C10").Value = 1.308019638061523
$ws.Range("E10").Value = 1552.907806659918
$ws.Range("F10").Value = 0.07776676636420352
$ws.Range("G10").Value = 0.0617006673332
$ws.Range("H10").Value = 0.05459836405230504
$ws.Range("I10").Value = 0.04753212285016289
$ws.Range("J10").Value = 0.04341415186756343
$ws.Range("K10").Value = 0.04034227403060377
$ws.Range("L10").Value = 0.03751456731739164
$ws.Range("M10").Value = 0.03560108334621516
$ws.Range("N10").Value = 0.03455703449140064
$ws.Range("O10").Value = 0.03367059286978676
$ws.Range("P10").Value = 0.03288617722939217
$ws.Range("Q10").Value = 0.03198171722661543
$ws.Range("R10").Value = 0.0315594877949309
$ws.Range("S10").Value = 0.03119102340018483
$ws.Range("T10").Value = 0.03084944930583816
$ws.Range("U10").Value = 0.03067011246826106
$ws.Range("V10").Value = 0.03050959950222781
$ws.Range("W10").Value = 0.03043229246611166
$ws.Range("X10").Value = 0.03030652034271507
$ws.Range("Y10").Value = 0.03027110734229859

# Row 11
$ws.Range("C11").Value = 1.129971504211426
$ws.Range("E11").Value = 1550.945829316684
$ws.Range("F11").Value = 0.07674030626819238
$ws.Range("G11").Value = 0.05976008487428367
$ws.Range("H11").Value = 0.05089620939222071
$ws.Range("I11").Value = 0.04387812362537052
$ws.Range("J11").Value = 0.04189585322602161
$ws.Range("K11").Value = 0.03883007158607752
$ws.Range("L11").Value = 0.03711183376213075
$ws.Range("M11").Value = 0.03495300451456498
$ws.Range("N11").Value = 0.03384983847175611
$ws.Range("O11").Value = 0.03290591735111401
$ws.Range("P11").Value = 0.03208196935223768
$ws.Range("Q11").Value = 0.03170988004105747
$ws.Range("R11").Value = 0.03162327162554464
$ws.Range("S11").Value = 0.03121733820735701
$ws.Range("T11").Value = 0.03101385456497073
$ws.Range("U11").Value = 0.03073330398648413
$ws.Range("V11").Value = 0.03062726213866102
$ws.Range("W11").Value = 0.03042808330854653
$ws.Range("X11").Value = 0.03031844216872494
$ws.Range("Y11").Value = 0.03023286216991587
